# Update cryptocurrency price/volume data to reflect the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text so Excel does not
# "helpfully" reinterpret values like 1.001 or 238.33 as numbers (which would
# mangle exact formatting such as trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.730.56"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "1.742.68"
$ws.Range("E3").Value = "  -4.98%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "238.33"
$ws.Range("E5").Value = "  -8.38%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "0.5042"
$ws.Range("E7").Value = "  -5.56%  "
$ws.Range("D8").Value = "41.83"
$ws.Range("E8").Value = "  -6.83%  "
$ws.Range("D9").Value = "0.2643"
$ws.Range("E9").Value = "  -12.40%  "
$ws.Range("D10").Value = "0.06135"
$ws.Range("E10").Value = "  -10.50%  "
$ws.Range("D11").Value = "1.745.63"
$ws.Range("E11").Value = "  -4.90%  "
$ws.Range("D12").Value = "0.06935"
$ws.Range("E12").Value = "  -4.57%  "
$ws.Range("E13").Value = "  -12.95%  "
$ws.Range("E14").Value = "  -9.36%  "
$ws.Range("D15").Value = "0.5936"
$ws.Range("E15").Value = "  -19.41%  "
$ws.Range("D16").Value = "76.53"
$ws.Range("E16").Value = "  -13.87%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "25.741.87"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").Value = "11.62"
$ws.Range("E20").Value = "  -15.93%  "
$ws.Range("D21").Value = "0.000006770"
$ws.Range("E21").Value = "  -13.81%  "
$ws.Range("D22").Value = "1.966.84"
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("D23").Value = "4.040"
$ws.Range("E23").Value = "  -11.30%  "
$ws.Range("D24").Value = "8.146"
$ws.Range("E24").Value = "  -11.33%  "
$ws.Range("D25").Value = "5.151"
$ws.Range("E25").Value = "  -13.39%  "
$ws.Range("D26").Value = "137.97"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").Value = "1.531"
$ws.Range("E27").Value = "  -8.93%  "
$ws.Range("D28").Value = "1.813"
$ws.Range("E28").Value = "  -17.03%  "
$ws.Range("D29").Value = "14.96"
$ws.Range("E29").Value = "  -11.49%  "
$ws.Range("D30").Value = "103.08"
$ws.Range("E30").Value = "  -6.57%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08107"
$ws.Range("E31").Value = "  -7.64%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "3.745"
$ws.Range("E32").Value = "  -11.13%  "
$ws.Range("D33").Value = "3.447"
$ws.Range("E33").Value = "  -13.59%  "
$ws.Range("D34").Value = "0.04498"
$ws.Range("E34").Value = "  -5.94%  "
$ws.Range("D35").Value = "0.9988"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "2.650"
$ws.Range("E36").Value = "  -9.43%  "
$ws.Range("D37").Value = "0.9793"
$ws.Range("E37").Value = "  -13.06%  "
$ws.Range("D38").Value = "0.6093"
$ws.Range("E38").Value = "  -16.20%  "
$ws.Range("D39").Value = "2.656"
$ws.Range("D40").Value = "0.01547"
$ws.Range("E40").Value = "  -8.97%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "0.9995"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "1.904"
$ws.Range("E42").Value = "  -16.13%  "
$ws.Range("D43").Value = "103.35"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").Value = "0.3790"
$ws.Range("E44").Value = "  -19.43%  "
$ws.Range("D45").Value = "5.107"
$ws.Range("E45").Value = "  -12.80%  "
$ws.Range("D46").Value = "0.7309"
$ws.Range("E46").Value = "  -18.99%  "
$ws.Range("D47").Value = "0.05339"
$ws.Range("E47").Value = "  -7.85%  "
$ws.Range("D48").Value = "0.1111"
$ws.Range("E48").Value = "  -9.45%  "
$ws.Range("D49").Value = "30.09"
$ws.Range("E49").Value = "  -13.09%  "
$ws.Range("D50").Value = "5.868"
$ws.Range("E50").Value = "  -19.98%  "
$ws.Range("D51").Value = "52.44"
$ws.Range("E51").Value = "  -12.43%  "
